$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value2 = 1.039483
$ws.Range("H2").Value2 = 3.118449
$ws.Range("I2").Value2 = 0.01443400247257094
$ws.Range("J2").Value2 = 0.01491631751620889
$ws.Range("M2").Value2 = 166.3936563333333
$ws.Range("N2").Value2 = 499.180969
$ws.Range("O2").Value2 = 0.6959913618211631
$ws.Range("P2").Value2 = 0.7009944564025758
$ws.Range("Q2").Value2 = 172.9633770663423
$ws.Range("R2").Value2 = 1556.670393597081
$ws.Range("S2").Value2 = 0.01004594103741469
$ws.Range("T2").Value2 = 0.01045625588880307

# Row 3
$ws.Range("G3").Value2 = 1.039483
$ws.Range("H3").Value2 = 3.118449
$ws.Range("I3").Value2 = 0.01443400247257094
$ws.Range("J3").Value2 = 0.01491631751620889
$ws.Range("O3").Value2 = 0.2039972194837954
$ws.Range("P3").Value2 = 0.2054636419703505
$ws.Range("Q3").Value2 = 50.696100454085
$ws.Range("R3").Value2 = 456.264904086765
$ws.Range("S3").Value2 = 0.002944496370426701
$ws.Range("T3").Value2 = 0.003064760921666411

# Row 4
$ws.Range("G4").Value2 = 1.039483
$ws.Range("H4").Value2 = 3.118449
$ws.Range("I4").Value2 = 0.01443400247257094
$ws.Range("J4").Value2 = 0.01491631751620889
$ws.Range("M4").Value2 = 7.402863
$ws.Range("N4").Value2 = 22.208589
$ws.Range("O4").Value2 = 0.03096469429353687
$ws.Range("P4").Value2 = 0.03118728224898178
$ws.Range("Q4").Value2 = 7.695150239828999
$ws.Range("R4").Value2 = 69.256352158461
$ws.Range("S4").Value2 = 0.0004469444739953145
$ws.Range("T4").Value2 = 0.0004651994044934376

# Row 5
$ws.Range("G5").Value2 = 1.039483
$ws.Range("H5").Value2 = 3.118449
$ws.Range("I5").Value2 = 0.01443400247257094
$ws.Range("J5").Value2 = 0.01491631751620889
$ws.Range("M5").Value2 = 5.118919500000001
$ws.Range("N5").Value2 = 10.237839
$ws.Range("O5").Value2 = 0.02141141574965316
$ws.Range("P5").Value2 = 0.0143768870013594
$ws.Range("Q5").Value2 = 5.3210297986185
$ws.Range("R5").Value2 = 31.926178791711
$ws.Range("S5").Value2 = 0.0003090524278717381
$ws.Range("T5").Value2 = 0.0002144502114069332

# Row 6
$ws.Range("G6").Value2 = 1.039483
$ws.Range("H6").Value2 = 3.118449
$ws.Range("I6").Value2 = 0.01443400247257094
$ws.Range("J6").Value2 = 0.01491631751620889
$ws.Range("M6").Value2 = 11.38837866666667
$ws.Range("N6").Value2 = 34.165136
$ws.Range("O6").Value2 = 0.04763530865185137
$ws.Range("P6").Value2 = 0.04797773237673265
$ws.Range("Q6").Value2 = 11.83802602156267
$ws.Range("R6").Value2 = 106.542234194064
$ws.Range("S6").Value2 = 0.0006875681628625027
$ws.Range("T6").Value2 = 0.0007156510898390395

# Row 7
$ws.Range("I7").Value2 = 0.8791289547788569
$ws.Range("J7").Value2 = 0.9085052224491242
$ws.Range("M7").Value2 = 166.3936563333333
$ws.Range("N7").Value2 = 499.180969
$ws.Range("O7").Value2 = 0.6959913618211631
$ws.Range("P7").Value2 = 0.7009944564025758
$ws.Range("Q7").Value2 = 10534.64644919593
$ws.Range("R7").Value2 = 94811.81804276342
$ws.Range("S7").Value2 = 0.6118661584529522
$ws.Range("T7").Value2 = 0.6368571245496251

# Row 8
$ws.Range("I8").Value2 = 0.8791289547788569
$ws.Range("J8").Value2 = 0.9085052224491242
$ws.Range("O8").Value2 = 0.2039972194837954
$ws.Range("P8").Value2 = 0.2054636419703505
$ws.Range("S8").Value2 = 0.1793398623425821
$ws.Range("T8").Value2 = 0.1866647917534805

# Row 9
$ws.Range("I9").Value2 = 0.8791289547788569
$ws.Range("J9").Value2 = 0.9085052224491242
$ws.Range("M9").Value2 = 7.402863
$ws.Range("N9").Value2 = 22.208589
$ws.Range("O9").Value2 = 0.03096469429353687
$ws.Range("P9").Value2 = 0.03118728224898178
$ws.Range("Q9").Value2 = 468.687004873581
$ws.Range("R9").Value2 = 4218.183043862229
$ws.Range("S9").Value2 = 0.0272219593293239
$ws.Range("T9").Value2 = 0.02833380879719482

# Row 10
$ws.Range("I10").Value2 = 0.8791289547788569
$ws.Range("J10").Value2 = 0.9085052224491242
$ws.Range("M10").Value2 = 5.118919500000001
$ws.Range("N10").Value2 = 10.237839
$ws.Range("O10").Value2 = 0.02141141574965316
$ws.Range("P10").Value2 = 0.0143768870013594
$ws.Range("Q10").Value2 = 324.0869172702465
$ws.Range("R10").Value2 = 1944.521503621479
$ws.Range("S10").Value2 = 0.01882339554832813
$ws.Range("T10").Value2 = 0.01306147692329595

# Row 11
$ws.Range("I11").Value2 = 0.8791289547788569
$ws.Range("J11").Value2 = 0.9085052224491242
$ws.Range("M11").Value2 = 11.38837866666667
$ws.Range("N11").Value2 = 34.165136
$ws.Range("O11").Value2 = 0.04763530865185137
$ws.Range("P11").Value2 = 0.04797773237673265
$ws.Range("Q11").Value2 = 721.0163267436108
$ws.Range("R11").Value2 = 6489.146940692497
$ws.Range("S11").Value2 = 0.04187757910567033
$ws.Range("T11").Value2 = 0.04358802042552804

# Row 12
$ws.Range("G12").Value2 = 0.3690693333333333
$ws.Range("H12").Value2 = 1.107208
$ws.Range("I12").Value2 = 0.005124804994293743
$ws.Range("J12").Value2 = 0.005296051365434103
$ws.Range("M12").Value2 = 166.3936563333333
$ws.Range("N12").Value2 = 499.180969
$ws.Range("O12").Value2 = 0.6959913618211631
$ws.Range("P12").Value2 = 0.7009944564025758
$ws.Range("Q12").Value2 = 61.4107958138391
$ws.Range("R12").Value2 = 552.697162324552
$ws.Range("S12").Value2 = 0.0035668200070464
$ws.Range("T12").Value2 = 0.003712502647992598

# Row 13
$ws.Range("G13").Value2 = 0.3690693333333333
$ws.Range("H13").Value2 = 1.107208
$ws.Range("I13").Value2 = 0.005124804994293743
$ws.Range("J13").Value2 = 0.005296051365434103
$ws.Range("O13").Value2 = 0.2039972194837954
$ws.Range("P13").Value2 = 0.2054636419703505
$ws.Range("Q13").Value2 = 17.99969407598667
$ws.Range("R13").Value2 = 161.99724668388
$ws.Range("S13").Value2 = 0.001045445969232592
$ws.Range("T13").Value2 = 0.001088146001604138

# Row 14
$ws.Range("G14").Value2 = 0.3690693333333333
$ws.Range("H14").Value2 = 1.107208
$ws.Range("I14").Value2 = 0.005124804994293743
$ws.Range("J14").Value2 = 0.005296051365434103
$ws.Range("M14").Value2 = 7.402863
$ws.Range("N14").Value2 = 22.208589
$ws.Range("O14").Value2 = 0.03096469429353687
$ws.Range("P14").Value2 = 0.03118728224898178
$ws.Range("Q14").Value2 = 2.732169712168
$ws.Range("R14").Value2 = 24.589527409512
$ws.Range("S14").Value2 = 0.0001586880199622967
$ws.Range("T14").Value2 = 0.0001651694487388987

# Row 15
$ws.Range("G15").Value2 = 0.3690693333333333
$ws.Range("H15").Value2 = 1.107208
$ws.Range("I15").Value2 = 0.005124804994293743
$ws.Range("J15").Value2 = 0.005296051365434103
$ws.Range("M15").Value2 = 5.118919500000001
$ws.Range("N15").Value2 = 10.237839
$ws.Range("O15").Value2 = 0.02141141574965316
$ws.Range("P15").Value2 = 0.0143768870013594
$ws.Range("Q15").Value2 = 1.889236207252
$ws.Range("R15").Value2 = 11.335417243512
$ws.Range("S15").Value2 = 0.0001097293303687222
$ws.Range("T15").Value2 = 0.00007614073203424127

# Row 16
$ws.Range("G16").Value2 = 0.3690693333333333
$ws.Range("H16").Value2 = 1.107208
$ws.Range("I16").Value2 = 0.005124804994293743
$ws.Range("J16").Value2 = 0.005296051365434103
$ws.Range("M16").Value2 = 11.38837866666667
$ws.Range("N16").Value2 = 34.165136
$ws.Range("O16").Value2 = 0.04763530865185137
$ws.Range("P16").Value2 = 0.04797773237673265
$ws.Range("Q16").Value2 = 4.203101322254223
$ws.Range("R16").Value2 = 37.827911900288
$ws.Range("S16").Value2 = 0.0002441216676837318
$ws.Range("T16").Value2 = 0.0002540925350642269

# Row 17
$ws.Range("G17").Value2 = 6.985879
$ws.Range("H17").Value2 = 13.971758
$ws.Range("I17").Value2 = 0.0970041787687547
$ws.Range("J17").Value2 = 0.06683039504177611
$ws.Range("M17").Value2 = 166.3936563333333
$ws.Range("N17").Value2 = 499.180969
$ws.Range("O17").Value2 = 0.6959913618211631
$ws.Range("P17").Value2 = 0.7009944564025758
$ws.Range("Q17").Value2 = 1162.40594951225
$ws.Range("R17").Value2 = 6974.435697073502
$ws.Range("S17").Value2 = 0.06751407048360913
$ws.Range("T17").Value2 = 0.04684773644347925

# Row 18
$ws.Range("G18").Value2 = 6.985879
$ws.Range("H18").Value2 = 13.971758
$ws.Range("I18").Value2 = 0.0970041787687547
$ws.Range("J18").Value2 = 0.06683039504177611
$ws.Range("O18").Value2 = 0.2039972194837954
$ws.Range("P18").Value2 = 0.2054636419703505
$ws.Range("Q18").Value2 = 340.704776840105
$ws.Range("R18").Value2 = 2044.22866104063
$ws.Range("S18").Value2 = 0.01978858274713498
$ws.Range("T18").Value2 = 0.01373121635960057

# Row 19
$ws.Range("G19").Value2 = 6.985879
$ws.Range("H19").Value2 = 13.971758
$ws.Range("I19").Value2 = 0.0970041787687547
$ws.Range("J19").Value2 = 0.06683039504177611
$ws.Range("M19").Value2 = 7.402863
$ws.Range("N19").Value2 = 22.208589
$ws.Range("O19").Value2 = 0.03096469429353687
$ws.Range("P19").Value2 = 0.03118728224898178
$ws.Range("Q19").Value2 = 51.715505171577
$ws.Range("R19").Value2 = 310.293031029462
$ws.Range("S19").Value2 = 0.003003704740770088
$ws.Range("T19").Value2 = 0.002084258392978824

# Row 20
$ws.Range("G20").Value2 = 6.985879
$ws.Range("H20").Value2 = 13.971758
$ws.Range("I20").Value2 = 0.0970041787687547
$ws.Range("J20").Value2 = 0.06683039504177611
$ws.Range("M20").Value2 = 5.118919500000001
$ws.Range("N20").Value2 = 10.237839
$ws.Range("O20").Value2 = 0.02141141574965316
$ws.Range("P20").Value2 = 0.0143768870013594
$ws.Range("Q20").Value2 = 35.7601522377405
$ws.Range("R20").Value2 = 143.040608950962
$ws.Range("S20").Value2 = 0.002076996801071485
$ws.Range("T20").Value2 = 0.0009608130377718249

# Row 21
$ws.Range("G21").Value2 = 6.985879
$ws.Range("H21").Value2 = 13.971758
$ws.Range("I21").Value2 = 0.0970041787687547
$ws.Range("J21").Value2 = 0.06683039504177611
$ws.Range("M21").Value2 = 11.38837866666667
$ws.Range("N21").Value2 = 34.165136
$ws.Range("O21").Value2 = 0.04763530865185137
$ws.Range("P21").Value2 = 0.04797773237673265
$ws.Range("Q21").Value2 = 79.55783537151467
$ws.Range("R21").Value2 = 477.347012229088
$ws.Range("S21").Value2 = 0.004620823996168997
$ws.Range("T21").Value2 = 0.003206370807945655

# Row 22
$ws.Range("G22").Value2 = 0.3102503333333334
$ws.Range("H22").Value2 = 0.930751
$ws.Range("I22").Value2 = 0.004308058985523854
$ws.Range("J22").Value2 = 0.004452013627456771
$ws.Range("M22").Value2 = 166.3936563333333
$ws.Range("N22").Value2 = 499.180969
$ws.Range("O22").Value2 = 0.6959913618211631
$ws.Range("P22").Value2 = 0.7009944564025758
$ws.Range("Q22").Value2 = 51.62368734196878
$ws.Range("R22").Value2 = 464.613186077719
$ws.Range("S22").Value2 = 0.002998371840140646
$ws.Range("T22").Value2 = 0.003120836872675919

# Row 23
$ws.Range("G23").Value2 = 0.3102503333333334
$ws.Range("H23").Value2 = 0.930751
$ws.Range("I23").Value2 = 0.004308058985523854
$ws.Range("J23").Value2 = 0.004452013627456771
$ws.Range("O23").Value2 = 0.2039972194837954
$ws.Range("P23").Value2 = 0.2054636419703505
$ws.Range("Q23").Value2 = 15.13106233058167
$ws.Range("R23").Value2 = 136.179560975235
$ws.Range("S23").Value2 = 0.0008788320544190468
$ws.Range("T23").Value2 = 0.0009147269339988993

# Row 24
$ws.Range("G24").Value2 = 0.3102503333333334
$ws.Range("H24").Value2 = 0.930751
$ws.Range("I24").Value2 = 0.004308058985523854
$ws.Range("J24").Value2 = 0.004452013627456771
$ws.Range("M24").Value2 = 7.402863
$ws.Range("N24").Value2 = 22.208589
$ws.Range("O24").Value2 = 0.03096469429353687
$ws.Range("P24").Value2 = 0.03118728224898178
$ws.Range("Q24").Value2 = 2.296740713371
$ws.Range("R24").Value2 = 20.670666420339
$ws.Range("S24").Value2 = 0.0001333977294852707
$ws.Range("T24").Value2 = 0.0001388462055758075

# Row 25
$ws.Range("G25").Value2 = 0.3102503333333334
$ws.Range("H25").Value2 = 0.930751
$ws.Range("I25").Value2 = 0.004308058985523854
$ws.Range("J25").Value2 = 0.004452013627456771
$ws.Range("M25").Value2 = 5.118919500000001
$ws.Range("N25").Value2 = 10.237839
$ws.Range("O25").Value2 = 0.02141141574965316
$ws.Range("P25").Value2 = 0.0143768870013594
$ws.Range("Q25").Value2 = 1.5881464811815
$ws.Range("R25").Value2 = 9.528878887089
$ws.Range("S25").Value2 = 0.00009224164201308024
$ws.Range("T25").Value2 = 0.00006400609685045818

# Row 26
$ws.Range("G26").Value2 = 0.3102503333333334
$ws.Range("H26").Value2 = 0.930751
$ws.Range("I26").Value2 = 0.004308058985523854
$ws.Range("J26").Value2 = 0.004452013627456771
$ws.Range("M26").Value2 = 11.38837866666667
$ws.Range("N26").Value2 = 34.165136
$ws.Range("O26").Value2 = 0.04763530865185137
$ws.Range("P26").Value2 = 0.04797773237673265
$ws.Range("Q26").Value2 = 3.533248277459556
$ws.Range("R26").Value2 = 31.799234497136
$ws.Range("S26").Value2 = 0.0002052157194658105
$ws.Range("T26").Value2 = 0.0002135975183556877

